$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the date (column A) and station (column B) values for rows 58-63,
# leaving the numeric columns (C:F) empty but keeping their existing style.
for ($r = 58; $r -le 63; $r++) {
    $ws.Range("A$r`:B$r").Clear()
}

# Update the active selection to match the diff (D62).
$ws.Range("D62").Select()
